$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 3.973

$ws.Range("G15").Value = 3.969

$ws.Range("G16").Value = 3.969

$ws.Range("G18").Value = 3.963

$ws.Range("G19").Value = 3.961

$ws.Range("B20").Value = 230121
$ws.Range("C20").Value = "DE MEL D.J."
$ws.Range("D20").Value = 3.957
$ws.Range("E20").Value = 4
$ws.Range("F20").Value = 3.921
$ws.Range("G20").Value = 3.96

$ws.Range("B21").Value = 230469
$ws.Range("C21").Value = "PEIRIS E.A.S.S."
$ws.Range("D21").Value = 4
$ws.Range("E21").Value = 3.96
$ws.Range("F21").Value = 3.937

$ws.Range("E23").Value = 3.921
$ws.Range("G23").Value = 3.953

$ws.Range("E24").Value = 3.96
$ws.Range("G24").Value = 3.953

$ws.Range("G25").Value = 3.951

$ws.Range("G28").Value = 3.942

$ws.Range("G29").Value = 3.942

$ws.Range("E31").Value = 3.934
$ws.Range("G31").Value = 3.936

$ws.Range("G32").Value = 3.935

$ws.Range("E34").Value = 3.96
$ws.Range("G34").Value = 3.921

$ws.Range("G35").Value = 3.92

$ws.Range("E36").Value = 3.96
$ws.Range("G36").Value = 3.916

$ws.Range("E37").Value = 3.908
$ws.Range("G37").Value = 3.911

$ws.Range("E38").Value = 3.934
$ws.Range("G38").Value = 3.909

$ws.Range("E39").Value = 3.939
$ws.Range("G39").Value = 3.901

$ws.Range("G40").Value = 3.893

$ws.Range("E41").Value = 3.869
$ws.Range("G41").Value = 3.891

$ws.Range("G42").Value = 3.881

$ws.Range("E43").Value = 3.96
$ws.Range("G43").Value = 3.881

$ws.Range("E44").Value = 3.882
$ws.Range("G44").Value = 3.876

$ws.Range("E45").Value = 3.947
$ws.Range("G45").Value = 3.876

$ws.Range("G46").Value = 3.875

$ws.Range("G47").Value = 3.871

$ws.Range("E48").Value = 3.817
$ws.Range("G48").Value = 3.866

$ws.Range("E49").Value = 3.947
$ws.Range("G49").Value = 3.865

$ws.Range("E50").Value = 3.908
$ws.Range("G50").Value = 3.855

$ws.Range("E51").Value = 3.96
$ws.Range("G51").Value = 3.855

$ws.Range("E52").Value = 3.895
$ws.Range("G52").Value = 3.853

$ws.Range("E53").Value = 3.908
$ws.Range("G53").Value = 3.851

$ws.Range("G54").Value = 3.851

$ws.Range("E55").Value = 3.869
$ws.Range("G55").Value = 3.841

$ws.Range("E56").Value = 3.908
$ws.Range("G56").Value = 3.836

$ws.Range("E57").Value = 3.908
$ws.Range("G57").Value = 3.831

$ws.Range("E58").Value = 3.817
$ws.Range("G58").Value = 3.83

$ws.Range("E59").Value = 3.908
$ws.Range("G59").Value = 3.83

$ws.Range("E61").Value = 3.947
$ws.Range("G61").Value = 3.821

$ws.Range("G62").Value = 3.818

$ws.Range("E63").Value = 3.96
$ws.Range("G63").Value = 3.813

$ws.Range("G64").Value = 3.808

$ws.Range("E65").Value = 3.908
$ws.Range("G65").Value = 3.796

$ws.Range("E66").Value = 3.869
$ws.Range("G66").Value = 3.796

$ws.Range("B67").Value = 230507
$ws.Range("C67").Value = "RAHMAN M.F.A."
$ws.Range("D67").Value = 3.857
$ws.Range("E67").Value = 3.877
$ws.Range("F67").Value = 3.645
$ws.Range("G67").Value = 3.796

$ws.Range("B68").Value = 230585
$ws.Range("C68").Value = "SARUKA U."
$ws.Range("D68").Value = 3.935
$ws.Range("E68").Value = 3.96
$ws.Range("F68").Value = 3.543
$ws.Range("G68").Value = 3.795

$ws.Range("B69").Value = 230726
$ws.Range("C69").Value = "WIJESINGHE U.G.S.K.D."
$ws.Range("D69").Value = 3.892
$ws.Range("E69").Value = 3.869
$ws.Range("F69").Value = 3.66
$ws.Range("G69").Value = 3.795

$ws.Range("B70").Value = 230016
$ws.Range("C70").Value = "ABISHEK L."
$ws.Range("D70").Value = 4
$ws.Range("E70").Value = 3.911
$ws.Range("F70").Value = 3.479
$ws.Range("G70").Value = 3.79

$ws.Range("B71").Value = 230070
$ws.Range("C71").Value = "BALASOORIYA B.R.B.D."
$ws.Range("D71").Value = 3.957
$ws.Range("E71").Value = 3.808
$ws.Range("F71").Value = 3.665
$ws.Range("G71").Value = 3.788

$ws.Range("E72").Value = 3.791
$ws.Range("G72").Value = 3.786

$ws.Range("G73").Value = 3.777

$ws.Range("E74").Value = 3.947
$ws.Range("G74").Value = 3.771

$ws.Range("E75").Value = 3.839
$ws.Range("G75").Value = 3.766

$ws.Range("E78").Value = 3.656
$ws.Range("G78").Value = 3.743

$ws.Range("E79").Value = 3.947
$ws.Range("G79").Value = 3.737

$ws.Range("E80").Value = 3.908
$ws.Range("G80").Value = 3.726

$ws.Range("E81").Value = 3.778
$ws.Range("G81").Value = 3.722

$ws.Range("E82").Value = 3.817
$ws.Range("G82").Value = 3.716

$ws.Range("E83").Value = 3.733
$ws.Range("G83").Value = 3.706

$ws.Range("E84").Value = 3.756
$ws.Range("G84").Value = 3.703

$ws.Range("E85").Value = 3.596
$ws.Range("G85").Value = 3.693

$ws.Range("E86").Value = 3.686
$ws.Range("G86").Value = 3.691

$ws.Range("E87").Value = 3.83
$ws.Range("G87").Value = 3.681

$ws.Range("E88").Value = 3.747
$ws.Range("G88").Value = 3.678

$ws.Range("E89").Value = 3.617
$ws.Range("G89").Value = 3.677

$ws.Range("E90").Value = 3.83
$ws.Range("G90").Value = 3.676

$ws.Range("E91").Value = 3.713
$ws.Range("G91").Value = 3.675

$ws.Range("E92").Value = 3.786
$ws.Range("G92").Value = 3.673

$ws.Range("E93").Value = 3.634
$ws.Range("G93").Value = 3.67

$ws.Range("E94").Value = 3.83
$ws.Range("G94").Value = 3.659

$ws.Range("E95").Value = 3.765
$ws.Range("G95").Value = 3.653

$ws.Range("E96").Value = 3.708
$ws.Range("G96").Value = 3.648

$ws.Range("E97").Value = 3.617
$ws.Range("G97").Value = 3.646

$ws.Range("E99").Value = 3.726
$ws.Range("G99").Value = 3.603

$ws.Range("E100").Value = 3.908
$ws.Range("G100").Value = 3.596

$ws.Range("B101").Value = 230259
$ws.Range("C101").Value = "IMBULPITIYA B.N."
$ws.Range("E101").Value = 3.551
$ws.Range("F101").Value = 3.447
$ws.Range("G101").Value = 3.583

$ws.Range("B102").Value = 230395
$ws.Range("C102").Value = "MANATUNGA K.D."
$ws.Range("E102").Value = 3.656
$ws.Range("F102").Value = 3.334
$ws.Range("G102").Value = 3.578

$ws.Range("E103").Value = 3.856
$ws.Range("G103").Value = 3.563

$ws.Range("E104").Value = 3.604
$ws.Range("G104").Value = 3.563

$ws.Range("E105").Value = 3.578
$ws.Range("G105").Value = 3.536

$ws.Range("E106").Value = 3.726
$ws.Range("G106").Value = 3.491

$ws.Range("E107").Value = 3.486
$ws.Range("G107").Value = 3.47

$ws.Range("B108").Value = 230581
$ws.Range("C108").Value = "SANTHOSH S."
$ws.Range("D108").Value = 3.792
$ws.Range("E108").Value = 3.629
$ws.Range("F108").Value = 3.014
$ws.Range("G108").Value = 3.465

$ws.Range("B109").Value = 230268
$ws.Range("C109").Value = "JAYAKODY J.A.C.P."
$ws.Range("D109").Value = 3.85
$ws.Range("E109").Value = 3.586
$ws.Range("F109").Value = 3.095
$ws.Range("G109").Value = 3.46

$ws.Range("E110").Value = 3.421
$ws.Range("G110").Value = 3.45

$ws.Range("E112").Value = 3.069
$ws.Range("G112").Value = 3.373

$ws.Range("E113").Value = 3.456
$ws.Range("G113").Value = 3.308

$ws.Range("E114").Value = 3.266
$ws.Range("G114").Value = 3.249

$ws.Range("E115").Value = 3.626
$ws.Range("G115").Value = 3.221

$ws.Range("F116").Value = 3.208
